$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Jeff Bamber's email address to add the third address (jeffrey.bamber@icr.ac.uk)
$ws.Range("E2").Value = "jeff.bamber@icr.ac.uk; jeff@icr.ac.uk; jeffrey.bamber@icr.ac.uk"

# Reflect the active cell selection moving to F2 (as seen after the edit/save)
$ws.Range("F2").Select()
